$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain numeric-looking text (e.g. "581.91",
# "11.08"). Force it to Text format first so Excel does not silently
# reinterpret the updated values as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "69.697.66"
$ws.Cells.Item(2, 5).Value = "  +2.04%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.383.33"
$ws.Cells.Item(3, 5).Value = "  +1.20%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.11%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "581.91"
$ws.Cells.Item(5, 5).Value = "  -0.48%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +1.19%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.10%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.595"
$ws.Cells.Item(8, 5).Value = "  +0.59%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +8.74%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.12%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "48.53"
$ws.Cells.Item(11, 5).Value = "  +0.59%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.0000285"
$ws.Cells.Item(12, 5).Value = "  +4.08%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "683.87"

# Row 14
$ws.Cells.Item(14, 4).Value = "8.63"
$ws.Cells.Item(14, 5).Value = "  +2.14%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.928.37"
$ws.Cells.Item(15, 5).Value = "  +1.16%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "69.709.30"
$ws.Cells.Item(16, 5).Value = "  +2.08%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.86%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.384.61"
$ws.Cells.Item(18, 5).Value = "  +1.30%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +1.24%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "11.28"
$ws.Cells.Item(20, 5).Value = "  +0.81%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.913"
$ws.Cells.Item(21, 5).Value = "  +1.97%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "17.35"
$ws.Cells.Item(22, 5).Value = "  +1.99%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -2.12%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "101.86"
$ws.Cells.Item(24, 5).Value = "  +1.70%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.65%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.09%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "9.73"
$ws.Cells.Item(27, 5).Value = "  +2.14%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +1.60%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "8.75"
$ws.Cells.Item(29, 5).Value = "  +2.42%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "6.92"
$ws.Cells.Item(30, 5).Value = "  -1.03%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "3.87"
$ws.Cells.Item(31, 5).Value = "  +16.30%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Bittensor"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(32, 4).Value = "557.40"
$ws.Cells.Item(32, 5).Value = "  -1.89%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Cosmos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(33, 4).Value = "11.08"
$ws.Cells.Item(33, 5).Value = "  +0.20%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  +0.41%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "57.91"
$ws.Cells.Item(35, 5).Value = "  +0.69%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.06%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "3.608.92"
$ws.Cells.Item(37, 5).Value = "  -2.67%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.30%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "35.42"
$ws.Cells.Item(39, 5).Value = "  +1.47%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +8.42%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +4.94%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "3.31"
$ws.Cells.Item(42, 5).Value = "  +3.64%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +3.28%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.25%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.25%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.11%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Mantle"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(47, 4).Value = "1.38"
$ws.Cells.Item(47, 5).Value = "  +3.87%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(48, 4).Value = "1.00"
$ws.Cells.Item(48, 5).Value = "  -0.08%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "130.45"
$ws.Cells.Item(49, 5).Value = "  -0.37%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "2.61"
$ws.Cells.Item(50, 5).Value = "  +1.23%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "7.48"
$ws.Cells.Item(51, 5).Value = "  +0.29%  "

